$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(100, 10, "2021-04-08 13:29:13.116392", 11, 9,  9,  10, 10, 10),
    @(30,  5,  "2021-04-08 13:29:35.558046", 5,  5,  5,  5,  5,  5),
    @(30,  5,  "2021-04-08 13:37:09.238480", 15, 15, 15, 1,  15, 15),
    @(40,  6,  "2021-04-08 13:39:48.093282", 25, 25, 25, 21, 25, 25),
    @(40,  6,  "2021-04-08 13:41:46.891262", 23, 23, 23, 20, 23, 23),
    @(40,  6,  "2021-04-08 13:42:09.620661", 12, 12, 12, 10, 12, 12),
    @(40,  6,  "2021-04-08 13:43:02.467024", 16, 16, 16, 6,  16, 16),
    @(40,  6,  "2021-04-08 13:43:20.036084", 20, 20, 20, 8,  20, 20),
    @(40,  6,  "2021-04-08 13:46:20.605873", 8,  8,  8,  7,  8,  8),
    @(40,  6,  "2021-04-08 13:54:29.560418", 22, 22, 22, 19, 22, 22)
)

$startRow = 19
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
    $ws.Cells.Item($row, 6).Value = $values[5]
    $ws.Cells.Item($row, 7).Value = $values[6]
    $ws.Cells.Item($row, 8).Value = $values[7]
    $ws.Cells.Item($row, 9).Value = $values[8]
}
